$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..4) {
    $ws.Cells.Item($row, 2).Value = "a-->b"
    $ws.Cells.Item($row, 4).Value = 1000000000
    $ws.Cells.Item($row, 5).Value = "a,b"
    $ws.Cells.Item($row, 6).Value = "1,1,1"
    $ws.Cells.Item($row, 7).Value = "a,b"
}
